$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "28.081.32"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "1.793.08"
$ws.Range("E3").Value = "  -0.41%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "316.71"
$ws.Range("E5").Value = "  +0.58%  "
$ws.Range("E6").Value = "  +0.00%  "
$ws.Range("D7").Value = "0.5392"
$ws.Range("E7").Value = "  -2.06%  "
$ws.Range("D8").Value = "0.3770"
$ws.Range("E8").Value = "  -2.11%  "
$ws.Range("D9").Value = "0.07438"
$ws.Range("E9").Value = "  -2.10%  "
$ws.Range("D10").Value = "41.71"
$ws.Range("E10").Value = "  -1.96%  "
$ws.Range("E11").Value = "  -3.10%  "
$ws.Range("D13").Value = "20.58"
$ws.Range("E13").Value = "  -2.76%  "
$ws.Range("D14").Value = "6.102"
$ws.Range("E14").Value = "  -1.33%  "
$ws.Range("D15").Value = "1.789.24"
$ws.Range("E15").Value = "  -0.94%  "
$ws.Range("D16").Value = "7.217"
$ws.Range("E16").Value = "  -2.63%  "
$ws.Range("D17").Value = "89.07"
$ws.Range("E17").Value = "  -3.13%  "
$ws.Range("E18").Value = "  -1.52%  "
$ws.Range("D19").Value = "0.06477"
$ws.Range("E19").Value = "  +0.45%  "
$ws.Range("D20").Value = "1.001"
$ws.Range("E20").Value = "  +0.04%  "
$ws.Range("D21").Value = "17.25"
$ws.Range("E21").Value = "  -0.65%  "
$ws.Range("D22").Value = "5.896"
$ws.Range("E22").Value = "  -1.26%  "
$ws.Range("D23").Value = "28.099.47"
$ws.Range("E23").Value = "  -1.11%  "
$ws.Range("D24").Value = "11.17"
$ws.Range("E24").Value = "  -2.31%  "
$ws.Range("D25").Value = "2.092"
$ws.Range("E25").Value = "  -1.60%  "
$ws.Range("D26").Value = "154.90"
$ws.Range("E26").Value = "  -2.43%  "
$ws.Range("D28").Value = "1.994.20"
$ws.Range("E28").Value = "  -1.00%  "
$ws.Range("D29").Value = "2.283"
$ws.Range("E29").Value = "  -5.11%  "
$ws.Range("D30").Value = "120.67"
$ws.Range("E30").Value = "  -2.52%  "
$ws.Range("D31").Value = "1.120"
$ws.Range("E31").Value = "  +0.06%  "
$ws.Range("D32").Value = "0.1055"
$ws.Range("E32").Value = "  +2.81%  "
$ws.Range("D33").Value = "3.656"
$ws.Range("E33").Value = "  -1.00%  "
$ws.Range("D34").Value = "5.544"
$ws.Range("E34").Value = "  -3.95%  "
$ws.Range("B35").Value = "Algorand"
$ws.Range("C35").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D35").Value = "0.2257"
$ws.Range("E35").Value = "  -2.21%  "
$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Value = "0.06507"
$ws.Range("E36").Value = "  +1.66%  "
$ws.Range("D37").Value = "0.02283"
$ws.Range("E37").Value = "  -1.90%  "
$ws.Range("D38").Value = "5.011"
$ws.Range("E38").Value = "  -2.70%  "
$ws.Range("D39").Value = "8.451"
$ws.Range("E39").Value = "  -3.71%  "
$ws.Range("D40").Value = "1.447"
$ws.Range("E40").Value = "  +4.48%  "
$ws.Range("D41").Value = "0.6161"
$ws.Range("E41").Value = "  -4.08%  "
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("E43").Value = "  +1.06%  "
$ws.Range("D44").Value = "1.000"
$ws.Range("E44").Value = "  -0.01%  "
$ws.Range("D45").Value = "13.28"
$ws.Range("E45").Value = "  -2.47%  "
$ws.Range("D46").Value = "3.671"
$ws.Range("E46").Value = "  -0.20%  "
$ws.Range("D47").Value = "0.5777"
$ws.Range("E47").Value = "  -3.32%  "
$ws.Range("D48").Value = "125.69"
$ws.Range("E48").Value = "  -0.65%  "
$ws.Range("D49").Value = "1.187"
$ws.Range("E49").Value = "  +3.14%  "
$ws.Range("D50").Value = "1.917"
$ws.Range("E50").Value = "  -3.23%  "
$ws.Range("D51").Value = "0.06811"
$ws.Range("E51").Value = "  -1.36%  "
